$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O should inherit the same formatting as column N (border, number
# format, font) since it's a continuation of the same table. Only touch
# rows 3-10 (rows 1-2 never had data/formatting in column N either).
$ws.Range("N3:N10").Copy() | Out-Null
$ws.Range("O3:O10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the new 2021 column of data.
$ws.Range("O4").Value = 2021
$ws.Range("O6").Value = 1860
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = 510
$ws.Range("O9").Value = 178
$ws.Range("O10").Value = 821

# Move the active selection to match the saved view state.
$ws.Range("P9").Select() | Out-Null
